$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 21, shifting existing rows 21-29 down to 22-30
$ws.Rows.Item(21).Insert()

# Populate the new row 21 with the new weekly data
$ws.Cells.Item(21, 1).Value = 8
$ws.Cells.Item(21, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(21, 3).Value = "Coquimbo"
$ws.Cells.Item(21, 4).Value = 44825
$ws.Cells.Item(21, 5).Value = 4
$ws.Cells.Item(21, 6).Value = 100112026
$ws.Cells.Item(21, 7).Value = "Haba"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 440
$ws.Cells.Item(21, 11).Value = 8000
$ws.Cells.Item(21, 12).Value = 9000
$ws.Cells.Item(21, 13).Value = 8500
$ws.Cells.Item(21, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(21, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(21, 16).Value = 340
$ws.Cells.Item(21, 17).Value = 25
$ws.Cells.Item(21, 18).Value = "Hortaliza"
